$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 0.3873
$ws.Cells.Item(2, 3).Value = 0.1322
$ws.Cells.Item(2, 4).Value = 0.2551
$ws.Cells.Item(3, 2).Value = 0.385
$ws.Cells.Item(3, 3).Value = 0.1307
$ws.Cells.Item(3, 4).Value = 0.2543
$ws.Cells.Item(4, 2).Value = 0.3821
$ws.Cells.Item(4, 3).Value = 0.1291
$ws.Cells.Item(4, 4).Value = 0.253
$ws.Cells.Item(5, 2).Value = 0.3789
$ws.Cells.Item(5, 3).Value = 0.1272
$ws.Cells.Item(5, 4).Value = 0.2517
$ws.Cells.Item(6, 2).Value = 0.3769
$ws.Cells.Item(6, 3).Value = 0.1261
$ws.Cells.Item(6, 4).Value = 0.2508
$ws.Cells.Item(7, 2).Value = 0.3747
$ws.Cells.Item(7, 3).Value = 0.1253
$ws.Cells.Item(7, 4).Value = 0.2494
$ws.Cells.Item(8, 2).Value = 0.3739
$ws.Cells.Item(8, 3).Value = 0.1255
$ws.Cells.Item(8, 4).Value = 0.2484
$ws.Cells.Item(9, 2).Value = 0.3737
$ws.Cells.Item(9, 3).Value = 0.1264
$ws.Cells.Item(9, 4).Value = 0.2473
$ws.Cells.Item(10, 2).Value = 0.3755
$ws.Cells.Item(10, 3).Value = 0.129
$ws.Cells.Item(10, 4).Value = 0.2466
$ws.Cells.Item(11, 2).Value = 0.3783
$ws.Cells.Item(11, 3).Value = 0.1324
$ws.Cells.Item(11, 4).Value = 0.2459
$ws.Cells.Item(12, 2).Value = 0.3816
$ws.Cells.Item(12, 3).Value = 0.1362
$ws.Cells.Item(12, 4).Value = 0.2454
$ws.Cells.Item(13, 2).Value = 0.3836
$ws.Cells.Item(13, 3).Value = 0.1391
$ws.Cells.Item(13, 4).Value = 0.2445
$ws.Cells.Item(14, 2).Value = 0.3846
$ws.Cells.Item(14, 3).Value = 0.1412
$ws.Cells.Item(14, 4).Value = 0.2434
$ws.Cells.Item(15, 2).Value = 0.3854
$ws.Cells.Item(15, 3).Value = 0.1431
$ws.Cells.Item(15, 4).Value = 0.2423
$ws.Cells.Item(16, 2).Value = 0.386
$ws.Cells.Item(16, 3).Value = 0.1447
$ws.Cells.Item(16, 4).Value = 0.2412
$ws.Cells.Item(17, 2).Value = 0.3871
$ws.Cells.Item(17, 3).Value = 0.1467
$ws.Cells.Item(17, 4).Value = 0.2404
$ws.Cells.Item(18, 2).Value = 0.3876
$ws.Cells.Item(18, 3).Value = 0.1479
$ws.Cells.Item(18, 4).Value = 0.2397
$ws.Cells.Item(19, 2).Value = 0.3885
$ws.Cells.Item(19, 3).Value = 0.1495
$ws.Cells.Item(19, 4).Value = 0.239
$ws.Cells.Item(20, 2).Value = 0.3894
$ws.Cells.Item(20, 3).Value = 0.1514
$ws.Cells.Item(20, 4).Value = 0.238
$ws.Cells.Item(21, 2).Value = 0.3899
$ws.Cells.Item(21, 3).Value = 0.1531
$ws.Cells.Item(21, 4).Value = 0.2368
$ws.Cells.Item(22, 2).Value = 0.3905
$ws.Cells.Item(22, 3).Value = 0.1552
$ws.Cells.Item(22, 4).Value = 0.2353
$ws.Cells.Item(23, 2).Value = 0.3908
$ws.Cells.Item(23, 3).Value = 0.1572
$ws.Cells.Item(23, 4).Value = 0.2336
$ws.Cells.Item(24, 2).Value = 0.3915
$ws.Cells.Item(24, 3).Value = 0.1598
$ws.Cells.Item(24, 4).Value = 0.2317
$ws.Cells.Item(25, 2).Value = 0.3919
$ws.Cells.Item(25, 3).Value = 0.1621
$ws.Cells.Item(25, 4).Value = 0.2298
$ws.Cells.Item(26, 2).Value = 0.3913
$ws.Cells.Item(26, 3).Value = 0.1638
$ws.Cells.Item(26, 4).Value = 0.2275
$ws.Cells.Item(27, 2).Value = 0.3906
$ws.Cells.Item(27, 3).Value = 0.1655
$ws.Cells.Item(27, 4).Value = 0.225
$ws.Cells.Item(28, 2).Value = 0.3886
$ws.Cells.Item(28, 3).Value = 0.1663
$ws.Cells.Item(28, 4).Value = 0.2223
$ws.Cells.Item(29, 2).Value = 0.3865
$ws.Cells.Item(29, 3).Value = 0.1669
$ws.Cells.Item(29, 4).Value = 0.2196
$ws.Cells.Item(30, 2).Value = 0.3828
$ws.Cells.Item(30, 3).Value = 0.1662
$ws.Cells.Item(30, 4).Value = 0.2166
$ws.Cells.Item(31, 2).Value = 0.3773
$ws.Cells.Item(31, 3).Value = 0.164
$ws.Cells.Item(31, 4).Value = 0.2133
$ws.Cells.Item(32, 2).Value = 0.3705
$ws.Cells.Item(32, 3).Value = 0.1609
$ws.Cells.Item(32, 4).Value = 0.2097
$ws.Cells.Item(33, 2).Value = 0.365
$ws.Cells.Item(33, 3).Value = 0.1587
$ws.Cells.Item(33, 4).Value = 0.2064
$ws.Cells.Item(34, 2).Value = 0.3602
$ws.Cells.Item(34, 3).Value = 0.1568
$ws.Cells.Item(34, 4).Value = 0.2034
$ws.Cells.Item(35, 2).Value = 0.3553
$ws.Cells.Item(35, 3).Value = 0.155
$ws.Cells.Item(35, 4).Value = 0.2003
$ws.Cells.Item(36, 2).Value = 0.3506
$ws.Cells.Item(36, 3).Value = 0.1535
$ws.Cells.Item(36, 4).Value = 0.1971
$ws.Cells.Item(37, 2).Value = 0.3461
$ws.Cells.Item(37, 3).Value = 0.1523
$ws.Cells.Item(37, 4).Value = 0.1937
$ws.Cells.Item(38, 2).Value = 0.3418
$ws.Cells.Item(38, 3).Value = 0.1509
$ws.Cells.Item(38, 4).Value = 0.1909
$ws.Cells.Item(39, 2).Value = 0.3374
$ws.Cells.Item(39, 3).Value = 0.1494
$ws.Cells.Item(39, 4).Value = 0.188
$ws.Cells.Item(40, 2).Value = 0.3319
$ws.Cells.Item(40, 3).Value = 0.1471
$ws.Cells.Item(40, 4).Value = 0.1848
$ws.Cells.Item(41, 2).Value = 0.3281
$ws.Cells.Item(41, 3).Value = 0.146
$ws.Cells.Item(41, 4).Value = 0.1821
$ws.Cells.Item(42, 2).Value = 0.3252
$ws.Cells.Item(42, 3).Value = 0.1456
$ws.Cells.Item(42, 4).Value = 0.1796
$ws.Cells.Item(43, 2).Value = 0.3226
$ws.Cells.Item(43, 3).Value = 0.145
$ws.Cells.Item(43, 4).Value = 0.1776
$ws.Cells.Item(44, 2).Value = 0.3201
$ws.Cells.Item(44, 3).Value = 0.1446
$ws.Cells.Item(44, 4).Value = 0.1755
$ws.Cells.Item(45, 2).Value = 0.3189
$ws.Cells.Item(45, 3).Value = 0.1447
$ws.Cells.Item(45, 4).Value = 0.1742
$ws.Cells.Item(46, 2).Value = 0.3172
$ws.Cells.Item(46, 3).Value = 0.1444
$ws.Cells.Item(46, 4).Value = 0.1728
$ws.Cells.Item(47, 2).Value = 0.3161
$ws.Cells.Item(47, 3).Value = 0.1444
$ws.Cells.Item(47, 4).Value = 0.1718
$ws.Cells.Item(48, 2).Value = 0.3152
$ws.Cells.Item(48, 3).Value = 0.1447
$ws.Cells.Item(48, 4).Value = 0.1705
$ws.Cells.Item(49, 2).Value = 0.3141
$ws.Cells.Item(49, 3).Value = 0.1449
$ws.Cells.Item(49, 4).Value = 0.1692
$ws.Cells.Item(50, 2).Value = 0.3139
$ws.Cells.Item(50, 3).Value = 0.1458
$ws.Cells.Item(50, 4).Value = 0.1681
$ws.Cells.Item(51, 2).Value = 0.3138
$ws.Cells.Item(51, 3).Value = 0.1466
$ws.Cells.Item(51, 4).Value = 0.1672
$ws.Cells.Item(52, 2).Value = 0.3147
$ws.Cells.Item(52, 3).Value = 0.148
$ws.Cells.Item(52, 4).Value = 0.1667
$ws.Cells.Item(53, 2).Value = 0.3153
$ws.Cells.Item(53, 3).Value = 0.1494
$ws.Cells.Item(53, 4).Value = 0.1659
$ws.Cells.Item(54, 2).Value = 0.3152
$ws.Cells.Item(54, 3).Value = 0.1505
$ws.Cells.Item(54, 4).Value = 0.1648
$ws.Cells.Item(55, 2).Value = 0.3152
$ws.Cells.Item(55, 3).Value = 0.1516
$ws.Cells.Item(55, 4).Value = 0.1636
$ws.Cells.Item(56, 2).Value = 0.316
$ws.Cells.Item(56, 3).Value = 0.1531
$ws.Cells.Item(56, 4).Value = 0.1629
$ws.Cells.Item(57, 2).Value = 0.3166
$ws.Cells.Item(57, 3).Value = 0.1543
$ws.Cells.Item(57, 4).Value = 0.1623
$ws.Cells.Item(58, 2).Value = 0.3167
$ws.Cells.Item(58, 3).Value = 0.1557
$ws.Cells.Item(58, 4).Value = 0.161
$ws.Cells.Item(59, 2).Value = 0.3171
$ws.Cells.Item(59, 3).Value = 0.1574
$ws.Cells.Item(59, 4).Value = 0.1598
$ws.Cells.Item(60, 2).Value = 0.3187
$ws.Cells.Item(60, 3).Value = 0.1599
$ws.Cells.Item(60, 4).Value = 0.1588
$ws.Cells.Item(61, 2).Value = 0.3181
$ws.Cells.Item(61, 3).Value = 0.1608
$ws.Cells.Item(61, 4).Value = 0.1573
$ws.Cells.Item(62, 2).Value = 0.3186
$ws.Cells.Item(62, 3).Value = 0.1614
$ws.Cells.Item(62, 4).Value = 0.1572
$ws.Cells.Item(63, 2).Value = 0.3186
$ws.Cells.Item(63, 3).Value = 0.1619
$ws.Cells.Item(63, 4).Value = 0.1567
$ws.Cells.Item(64, 2).Value = 0.318
$ws.Cells.Item(64, 3).Value = 0.1618
$ws.Cells.Item(64, 4).Value = 0.1562
$ws.Cells.Item(65, 2).Value = 0.3174
$ws.Cells.Item(65, 3).Value = 0.1622
$ws.Cells.Item(65, 4).Value = 0.1552
$ws.Cells.Item(66, 2).Value = 0.3164
$ws.Cells.Item(66, 3).Value = 0.162
$ws.Cells.Item(66, 4).Value = 0.1544
$ws.Cells.Item(67, 2).Value = 0.3154
$ws.Cells.Item(67, 3).Value = 0.1617
$ws.Cells.Item(67, 4).Value = 0.1537
$ws.Cells.Item(68, 2).Value = 0.3139
$ws.Cells.Item(68, 3).Value = 0.1609
$ws.Cells.Item(68, 4).Value = 0.153
$ws.Cells.Item(69, 2).Value = 0.3121
$ws.Cells.Item(69, 3).Value = 0.1597
$ws.Cells.Item(69, 4).Value = 0.1524
$ws.Cells.Item(70, 2).Value = 0.3089
$ws.Cells.Item(70, 3).Value = 0.1573
$ws.Cells.Item(70, 4).Value = 0.1515
$ws.Cells.Item(71, 2).Value = 0.3065
$ws.Cells.Item(71, 3).Value = 0.1557
$ws.Cells.Item(71, 4).Value = 0.1508
$ws.Cells.Item(72, 2).Value = 0.3038
$ws.Cells.Item(72, 3).Value = 0.1538
$ws.Cells.Item(72, 4).Value = 0.1499
$ws.Cells.Item(73, 2).Value = 0.3023
$ws.Cells.Item(73, 3).Value = 0.1526
$ws.Cells.Item(73, 4).Value = 0.1497
$ws.Cells.Item(74, 2).Value = 0.3012
$ws.Cells.Item(74, 3).Value = 0.1511
$ws.Cells.Item(74, 4).Value = 0.1501
$ws.Cells.Item(75, 2).Value = 0.3007
$ws.Cells.Item(75, 3).Value = 0.1499
$ws.Cells.Item(75, 4).Value = 0.1508
$ws.Cells.Item(76, 2).Value = 0.2998
$ws.Cells.Item(76, 3).Value = 0.1484
$ws.Cells.Item(76, 4).Value = 0.1514
$ws.Cells.Item(77, 2).Value = 0.2987
$ws.Cells.Item(77, 3).Value = 0.147
$ws.Cells.Item(77, 4).Value = 0.1517
$ws.Cells.Item(78, 2).Value = 0.2976
$ws.Cells.Item(78, 3).Value = 0.1457
$ws.Cells.Item(78, 4).Value = 0.1519
$ws.Cells.Item(79, 2).Value = 0.2964
$ws.Cells.Item(79, 3).Value = 0.1441
$ws.Cells.Item(79, 4).Value = 0.1523
$ws.Cells.Item(80, 2).Value = 0.2956
$ws.Cells.Item(80, 3).Value = 0.1427
$ws.Cells.Item(80, 4).Value = 0.1528
$ws.Cells.Item(81, 2).Value = 0.2951
$ws.Cells.Item(81, 3).Value = 0.1415
$ws.Cells.Item(81, 4).Value = 0.1536
$ws.Cells.Item(82, 2).Value = 0.2933
$ws.Cells.Item(82, 3).Value = 0.1401
$ws.Cells.Item(82, 4).Value = 0.1531
$ws.Cells.Item(83, 2).Value = 0.2917
$ws.Cells.Item(83, 3).Value = 0.1386
$ws.Cells.Item(83, 4).Value = 0.1531
$ws.Cells.Item(84, 2).Value = 0.2906
$ws.Cells.Item(84, 3).Value = 0.1372
$ws.Cells.Item(84, 4).Value = 0.1534
$ws.Cells.Item(85, 2).Value = 0.2893
$ws.Cells.Item(85, 3).Value = 0.1357
$ws.Cells.Item(85, 4).Value = 0.1536
$ws.Cells.Item(86, 2).Value = 0.2888
$ws.Cells.Item(86, 3).Value = 0.1345
$ws.Cells.Item(86, 4).Value = 0.1543
$ws.Cells.Item(87, 2).Value = 0.2885
$ws.Cells.Item(87, 3).Value = 0.1336
$ws.Cells.Item(87, 4).Value = 0.1549
$ws.Cells.Item(88, 2).Value = 0.2893
$ws.Cells.Item(88, 3).Value = 0.133
$ws.Cells.Item(88, 4).Value = 0.1563
$ws.Cells.Item(89, 2).Value = 0.2902
$ws.Cells.Item(89, 3).Value = 0.1327
$ws.Cells.Item(89, 4).Value = 0.1575
$ws.Cells.Item(90, 2).Value = 0.2917
$ws.Cells.Item(90, 3).Value = 0.1328
$ws.Cells.Item(90, 4).Value = 0.1589
$ws.Cells.Item(91, 2).Value = 0.2929
$ws.Cells.Item(91, 3).Value = 0.1326
$ws.Cells.Item(91, 4).Value = 0.1603
$ws.Cells.Item(92, 2).Value = 0.2935
$ws.Cells.Item(92, 3).Value = 0.1321
$ws.Cells.Item(92, 4).Value = 0.1614
$ws.Cells.Item(93, 2).Value = 0.294
$ws.Cells.Item(93, 3).Value = 0.1322
$ws.Cells.Item(93, 4).Value = 0.1618
$ws.Cells.Item(94, 2).Value = 0.295
$ws.Cells.Item(94, 3).Value = 0.1328
$ws.Cells.Item(94, 4).Value = 0.1622
$ws.Cells.Item(95, 2).Value = 0.2947
$ws.Cells.Item(95, 3).Value = 0.1324
$ws.Cells.Item(95, 4).Value = 0.1623
$ws.Cells.Item(96, 2).Value = 0.2944
$ws.Cells.Item(96, 3).Value = 0.1324
$ws.Cells.Item(96, 4).Value = 0.162
$ws.Cells.Item(97, 2).Value = 0.2977
$ws.Cells.Item(97, 3).Value = 0.1343
$ws.Cells.Item(97, 4).Value = 0.1634
$ws.Cells.Item(98, 2).Value = 0.3012
$ws.Cells.Item(98, 3).Value = 0.1363
$ws.Cells.Item(98, 4).Value = 0.1649
$ws.Cells.Item(99, 2).Value = 0.3048
$ws.Cells.Item(99, 3).Value = 0.1385
$ws.Cells.Item(99, 4).Value = 0.1663
$ws.Cells.Item(100, 2).Value = 0.3079
$ws.Cells.Item(100, 3).Value = 0.1405
$ws.Cells.Item(100, 4).Value = 0.1674
$ws.Cells.Item(101, 2).Value = 0.3124
$ws.Cells.Item(101, 3).Value = 0.1436
$ws.Cells.Item(101, 4).Value = 0.1688
$ws.Cells.Item(102, 2).Value = 0.3166
$ws.Cells.Item(102, 3).Value = 0.1468
$ws.Cells.Item(102, 4).Value = 0.1698
$ws.Cells.Item(103, 2).Value = 0.3207
$ws.Cells.Item(103, 3).Value = 0.1501
$ws.Cells.Item(103, 4).Value = 0.1707
$ws.Cells.Item(104, 2).Value = 0.324
$ws.Cells.Item(104, 3).Value = 0.1528
$ws.Cells.Item(104, 4).Value = 0.1712
$ws.Cells.Item(105, 2).Value = 0.3257
$ws.Cells.Item(105, 3).Value = 0.1541
$ws.Cells.Item(105, 4).Value = 0.1716
$ws.Cells.Item(106, 2).Value = 0.3269
$ws.Cells.Item(106, 3).Value = 0.1552
$ws.Cells.Item(106, 4).Value = 0.1717
$ws.Cells.Item(107, 2).Value = 0.3279
$ws.Cells.Item(107, 3).Value = 0.1563
$ws.Cells.Item(107, 4).Value = 0.1716
$ws.Cells.Item(108, 2).Value = 0.3284
$ws.Cells.Item(108, 3).Value = 0.1577
$ws.Cells.Item(108, 4).Value = 0.1707
$ws.Cells.Item(109, 2).Value = 0.3293
$ws.Cells.Item(109, 3).Value = 0.1591
$ws.Cells.Item(109, 4).Value = 0.1702
$ws.Cells.Item(110, 2).Value = 0.3307
$ws.Cells.Item(110, 3).Value = 0.1604
$ws.Cells.Item(110, 4).Value = 0.1704
$ws.Cells.Item(111, 2).Value = 0.3324
$ws.Cells.Item(111, 3).Value = 0.1617
$ws.Cells.Item(111, 4).Value = 0.1706
$ws.Cells.Item(112, 2).Value = 0.3341
$ws.Cells.Item(112, 3).Value = 0.163
$ws.Cells.Item(112, 4).Value = 0.1711
$ws.Cells.Item(113, 2).Value = 0.3351
$ws.Cells.Item(113, 3).Value = 0.1633
$ws.Cells.Item(113, 4).Value = 0.1718
$ws.Cells.Item(114, 2).Value = 0.3352
$ws.Cells.Item(114, 3).Value = 0.1635
$ws.Cells.Item(114, 4).Value = 0.1717
$ws.Cells.Item(115, 2).Value = 0.3359
$ws.Cells.Item(115, 3).Value = 0.1642
$ws.Cells.Item(115, 4).Value = 0.1716
$ws.Cells.Item(116, 2).Value = 0.3365
$ws.Cells.Item(116, 3).Value = 0.1647
$ws.Cells.Item(116, 4).Value = 0.1718
$ws.Cells.Item(117, 2).Value = 0.3341
$ws.Cells.Item(117, 3).Value = 0.1638
$ws.Cells.Item(117, 4).Value = 0.1703
$ws.Cells.Item(118, 2).Value = 0.3316
$ws.Cells.Item(118, 3).Value = 0.163
$ws.Cells.Item(118, 4).Value = 0.1687
$ws.Cells.Item(119, 2).Value = 0.33
$ws.Cells.Item(119, 3).Value = 0.1628
$ws.Cells.Item(119, 4).Value = 0.1672
$ws.Cells.Item(120, 2).Value = 0.3281
$ws.Cells.Item(120, 3).Value = 0.1626
$ws.Cells.Item(120, 4).Value = 0.1655
$ws.Cells.Item(121, 2).Value = 0.3246
$ws.Cells.Item(121, 3).Value = 0.1611
$ws.Cells.Item(121, 4).Value = 0.1635
$ws.Cells.Item(122, 2).Value = 0.3208
$ws.Cells.Item(122, 3).Value = 0.1594
$ws.Cells.Item(122, 4).Value = 0.1614
$ws.Cells.Item(123, 2).Value = 0.3163
$ws.Cells.Item(123, 3).Value = 0.1574
$ws.Cells.Item(123, 4).Value = 0.1589
$ws.Cells.Item(124, 2).Value = 0.3138
$ws.Cells.Item(124, 3).Value = 0.1566
$ws.Cells.Item(124, 4).Value = 0.1572
$ws.Cells.Item(125, 2).Value = 0.3108
$ws.Cells.Item(125, 3).Value = 0.1558
$ws.Cells.Item(125, 4).Value = 0.1551
$ws.Cells.Item(126, 2).Value = 0.3091
$ws.Cells.Item(126, 3).Value = 0.1558
$ws.Cells.Item(126, 4).Value = 0.1533
$ws.Cells.Item(127, 2).Value = 0.3073
$ws.Cells.Item(127, 3).Value = 0.1561
$ws.Cells.Item(127, 4).Value = 0.1512
$ws.Cells.Item(128, 2).Value = 0.3052
$ws.Cells.Item(128, 3).Value = 0.1556
$ws.Cells.Item(128, 4).Value = 0.1496
$ws.Cells.Item(129, 2).Value = 0.3037
$ws.Cells.Item(129, 3).Value = 0.1551
$ws.Cells.Item(129, 4).Value = 0.1486
$ws.Cells.Item(130, 2).Value = 0.3027
$ws.Cells.Item(130, 3).Value = 0.1559
$ws.Cells.Item(130, 4).Value = 0.1468
$ws.Cells.Item(131, 2).Value = 0.3032
$ws.Cells.Item(131, 3).Value = 0.1581
$ws.Cells.Item(131, 4).Value = 0.1451
$ws.Cells.Item(132, 2).Value = 0.3045
$ws.Cells.Item(132, 3).Value = 0.1598
$ws.Cells.Item(132, 4).Value = 0.1447
$ws.Cells.Item(133, 2).Value = 0.3045
$ws.Cells.Item(133, 3).Value = 0.1604
$ws.Cells.Item(133, 4).Value = 0.1441
$ws.Cells.Item(134, 2).Value = 0.3046
$ws.Cells.Item(134, 3).Value = 0.1612
$ws.Cells.Item(134, 4).Value = 0.1435
$ws.Cells.Item(135, 2).Value = 0.3045
$ws.Cells.Item(135, 3).Value = 0.1617
$ws.Cells.Item(135, 4).Value = 0.1428
$ws.Cells.Item(136, 2).Value = 0.3042
$ws.Cells.Item(136, 3).Value = 0.1623
$ws.Cells.Item(136, 4).Value = 0.1419
$ws.Cells.Item(137, 2).Value = 0.3045
$ws.Cells.Item(137, 3).Value = 0.1632
$ws.Cells.Item(137, 4).Value = 0.1414
$ws.Cells.Item(138, 2).Value = 0.3049
$ws.Cells.Item(138, 3).Value = 0.1639
$ws.Cells.Item(138, 4).Value = 0.141
$ws.Cells.Item(139, 2).Value = 0.3036
$ws.Cells.Item(139, 3).Value = 0.1633
$ws.Cells.Item(139, 4).Value = 0.1403
$ws.Cells.Item(140, 2).Value = 0.3024
$ws.Cells.Item(140, 3).Value = 0.1624
$ws.Cells.Item(140, 4).Value = 0.14
$ws.Cells.Item(141, 2).Value = 0.3024
$ws.Cells.Item(141, 3).Value = 0.1629
$ws.Cells.Item(141, 4).Value = 0.1395
$ws.Cells.Item(142, 2).Value = 0.3012
$ws.Cells.Item(142, 3).Value = 0.1624
$ws.Cells.Item(142, 4).Value = 0.1389
$ws.Cells.Item(143, 2).Value = 0.3
$ws.Cells.Item(143, 3).Value = 0.1619
$ws.Cells.Item(143, 4).Value = 0.1381
$ws.Cells.Item(144, 2).Value = 0.2974
$ws.Cells.Item(144, 3).Value = 0.1608
$ws.Cells.Item(144, 4).Value = 0.1367
$ws.Cells.Item(145, 2).Value = 0.2961
$ws.Cells.Item(145, 3).Value = 0.1603
$ws.Cells.Item(145, 4).Value = 0.1358
$ws.Cells.Item(146, 2).Value = 0.295
$ws.Cells.Item(146, 3).Value = 0.1602
$ws.Cells.Item(146, 4).Value = 0.1348
$ws.Cells.Item(147, 2).Value = 0.2928
$ws.Cells.Item(147, 3).Value = 0.1587
$ws.Cells.Item(147, 4).Value = 0.1341
$ws.Cells.Item(148, 2).Value = 0.2914
$ws.Cells.Item(148, 3).Value = 0.1585
$ws.Cells.Item(148, 4).Value = 0.1328
$ws.Cells.Item(149, 2).Value = 0.2905
$ws.Cells.Item(149, 3).Value = 0.1596
$ws.Cells.Item(149, 4).Value = 0.1309
$ws.Cells.Item(150, 2).Value = 0.2874
$ws.Cells.Item(150, 3).Value = 0.1588
$ws.Cells.Item(150, 4).Value = 0.1286
$ws.Cells.Item(151, 2).Value = 0.2847
$ws.Cells.Item(151, 3).Value = 0.1581
$ws.Cells.Item(151, 4).Value = 0.1266
$ws.Cells.Item(152, 2).Value = 0.2806
$ws.Cells.Item(152, 3).Value = 0.1573
$ws.Cells.Item(152, 4).Value = 0.1233
$ws.Cells.Item(153, 2).Value = 0.2777
$ws.Cells.Item(153, 3).Value = 0.1577
$ws.Cells.Item(153, 4).Value = 0.12
$ws.Cells.Item(154, 2).Value = 0.2748
$ws.Cells.Item(154, 3).Value = 0.1581
$ws.Cells.Item(154, 4).Value = 0.1168
$ws.Cells.Item(155, 2).Value = 0.272
$ws.Cells.Item(155, 3).Value = 0.1581
$ws.Cells.Item(155, 4).Value = 0.1139
$ws.Cells.Item(156, 2).Value = 0.2695
$ws.Cells.Item(156, 3).Value = 0.1584
$ws.Cells.Item(156, 4).Value = 0.1111
$ws.Cells.Item(157, 2).Value = 0.2651
$ws.Cells.Item(157, 3).Value = 0.1571
$ws.Cells.Item(157, 4).Value = 0.1081
$ws.Cells.Item(158, 2).Value = 0.2611
$ws.Cells.Item(158, 3).Value = 0.156
$ws.Cells.Item(158, 4).Value = 0.1051
$ws.Cells.Item(159, 2).Value = 0.2575
$ws.Cells.Item(159, 3).Value = 0.1551
$ws.Cells.Item(159, 4).Value = 0.1024
$ws.Cells.Item(160, 2).Value = 0.2542
$ws.Cells.Item(160, 3).Value = 0.1546
$ws.Cells.Item(160, 4).Value = 0.0997
$ws.Cells.Item(161, 2).Value = 0.2508
$ws.Cells.Item(161, 3).Value = 0.1536
$ws.Cells.Item(161, 4).Value = 0.0973
$ws.Cells.Item(162, 2).Value = 0.2477
$ws.Cells.Item(162, 3).Value = 0.1524
$ws.Cells.Item(162, 4).Value = 0.0952
$ws.Cells.Item(163, 2).Value = 0.2457
$ws.Cells.Item(163, 3).Value = 0.152
$ws.Cells.Item(163, 4).Value = 0.0936
$ws.Cells.Item(164, 2).Value = 0.2436
$ws.Cells.Item(164, 3).Value = 0.1513
$ws.Cells.Item(164, 4).Value = 0.0923

$ws.Cells.Item(165, 1).NumberFormat = "@"
$ws.Cells.Item(165, 1).Value = "26-08-2021"
$ws.Cells.Item(165, 1).Style = "Normal"
$ws.Cells.Item(165, 2).Value = 0.2413
$ws.Cells.Item(165, 3).Value = 0.1503
$ws.Cells.Item(165, 4).Value = 0.091
$ws.Cells.Item(166, 1).NumberFormat = "@"
$ws.Cells.Item(166, 1).Value = "27-08-2021"
$ws.Cells.Item(166, 1).Style = "Normal"
$ws.Cells.Item(166, 2).Value = 0.2396
$ws.Cells.Item(166, 3).Value = 0.1497
$ws.Cells.Item(166, 4).Value = 0.0898
$ws.Cells.Item(167, 1).NumberFormat = "@"
$ws.Cells.Item(167, 1).Value = "30-08-2021"
$ws.Cells.Item(167, 1).Style = "Normal"
$ws.Cells.Item(167, 2).Value = 0.2371
$ws.Cells.Item(167, 3).Value = 0.1485
$ws.Cells.Item(167, 4).Value = 0.0886
$ws.Cells.Item(168, 1).NumberFormat = "@"
$ws.Cells.Item(168, 1).Value = "31-08-2021"
$ws.Cells.Item(168, 1).Style = "Normal"
$ws.Cells.Item(168, 2).Value = 0.2339
$ws.Cells.Item(168, 3).Value = 0.1465
$ws.Cells.Item(168, 4).Value = 0.0874
$ws.Cells.Item(169, 1).NumberFormat = "@"
$ws.Cells.Item(169, 1).Value = "01-09-2021"
$ws.Cells.Item(169, 1).Style = "Normal"
$ws.Cells.Item(169, 2).Value = 0.2294
$ws.Cells.Item(169, 3).Value = 0.1428
$ws.Cells.Item(169, 4).Value = 0.0866
$ws.Cells.Item(170, 1).NumberFormat = "@"
$ws.Cells.Item(170, 1).Value = "02-09-2021"
$ws.Cells.Item(170, 1).Style = "Normal"
$ws.Cells.Item(170, 2).Value = 0.2293
$ws.Cells.Item(170, 3).Value = 0.141
$ws.Cells.Item(170, 4).Value = 0.0883
$ws.Cells.Item(171, 1).NumberFormat = "@"
$ws.Cells.Item(171, 1).Value = "03-09-2021"
$ws.Cells.Item(171, 1).Style = "Normal"
$ws.Cells.Item(171, 2).Value = 0.2278
$ws.Cells.Item(171, 3).Value = 0.1381
$ws.Cells.Item(171, 4).Value = 0.0897
$ws.Cells.Item(172, 1).NumberFormat = "@"
$ws.Cells.Item(172, 1).Value = "06-09-2021"
$ws.Cells.Item(172, 1).Style = "Normal"
$ws.Cells.Item(172, 2).Value = 0.227
$ws.Cells.Item(172, 3).Value = 0.1359
$ws.Cells.Item(172, 4).Value = 0.0912
$ws.Cells.Item(173, 1).NumberFormat = "@"
$ws.Cells.Item(173, 1).Value = "07-09-2021"
$ws.Cells.Item(173, 1).Style = "Normal"
$ws.Cells.Item(173, 2).Value = 0.2259
$ws.Cells.Item(173, 3).Value = 0.1332
$ws.Cells.Item(173, 4).Value = 0.0926
$ws.Cells.Item(174, 1).NumberFormat = "@"
$ws.Cells.Item(174, 1).Value = "08-09-2021"
$ws.Cells.Item(174, 1).Style = "Normal"
$ws.Cells.Item(174, 2).Value = 0.225
$ws.Cells.Item(174, 3).Value = 0.1308
$ws.Cells.Item(174, 4).Value = 0.0942
